$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rows 65-70: add "Date Opened" (F) and "Status" (G) values ---
$ws.Range("F65").Value = 43677
$ws.Range("G65").Value = "DONE"

$ws.Range("F66").Value = 43677
$ws.Range("G66").Value = "DONE"

$ws.Range("F67").Value = 43677
$ws.Range("G67").Value = "DONE"

$ws.Range("F68").Value = 43677
$ws.Range("G68").Value = "OPEN"

$ws.Range("F69").Value = 43677
$ws.Range("G69").Value = "OPEN"

$ws.Range("F70").Value = 43677
$ws.Range("G70").Value = "OPEN"

# --- Row 71: same Date Opened / Status addition, plus the row grows a touch taller ---
$ws.Range("F71").Value = 43677
$ws.Range("G71").Value = "OPEN"
$ws.Rows.Item(71).RowHeight = 48

# --- New issue row 117 ---
$ws.Range("A117").Value = 117
$ws.Range("B117").Value = "EPMD segments do not line up properly when ratio is different"
$ws.Range("C117").Value = -1
$ws.Range("F117").Value = 43677
$ws.Range("G117").Value = "OPEN"

# Match the author's final selection in the frozen (bottom-left) pane
$ws.Range("A66:B66").Select() | Out-Null
